$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: true -> (blank, cell cleared)
$ws.Range("B7").ClearContents()

# Date: 2023-10-31 -> 2025-11-18
# Force the new value to stay a text string (not be auto-converted to a date
# serial number), then restore the original (unformatted) cell style so the
# cell keeps looking like every other "plain" metadata cell.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2025-11-18"
$ws.Range("B9").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
